$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename topic codes from "en7-uXX" to "en5-uXX" (flipbooks now target class 5 instead of 7)
$ws.Range("A2").Value = "en5-u01"
$ws.Range("A3").Value = "en5-u01"
$ws.Range("A4").Value = "en5-u02"
$ws.Range("A5").Value = "en5-u02"
